# Refresh the crypto price/volume table (columns D and E, rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $text into $range as a literal TEXT value (matching the sheet's existing
# string cells in column D) without letting Excel auto-coerce number-looking
# strings (e.g. "306.00") into numbers -- which would silently drop trailing
# zeros -- and without leaving a NumberFormat/style side effect behind.
function Set-TextCell($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
    $excel.CutCopyMode = $false
}

# New Price (column D, $null = unchanged) / Volume(1h) (column E) per row.
$updates = @(
    @{ Row = 2; D = '41.937.42'; E = '  -1.25%  ' },
    @{ Row = 3; D = '2.246.93'; E = '  -1.40%  ' },
    @{ Row = 4; D = $null; E = '  -0.09%  ' },
    @{ Row = 5; D = '306.00'; E = '  -0.45%  ' },
    @{ Row = 6; D = '96.29'; E = '  -1.15%  ' },
    @{ Row = 7; D = $null; E = '  -1.52%  ' },
    @{ Row = 8; D = $null; E = '  -0.03%  ' },
    @{ Row = 9; D = '0.488'; E = '  -0.82%  ' },
    @{ Row = 10; D = '34.57'; E = '  -3.88%  ' },
    @{ Row = 11; D = '0.0812'; E = '  +2.05%  ' },
    @{ Row = 12; D = $null; E = '  +0.91%  ' },
    @{ Row = 13; D = '6.75'; E = '  +0.93%  ' },
    @{ Row = 14; D = '2.593.97'; E = '  -1.56%  ' },
    @{ Row = 15; D = '14.34'; E = '  -0.54%  ' },
    @{ Row = 16; D = '2.234.21'; E = '  -2.47%  ' },
    @{ Row = 17; D = '0.779'; E = '  -2.35%  ' },
    @{ Row = 18; D = '41.830.06'; E = '  -1.29%  ' },
    @{ Row = 19; D = '12.12'; E = '  -3.17%  ' },
    @{ Row = 20; D = '0.0₃0899'; E = '  -1.17%  ' },
    @{ Row = 21; D = '5.91'; E = '  -0.84%  ' },
    @{ Row = 22; D = '67.11'; E = '  -0.96%  ' },
    @{ Row = 23; D = '235.12'; E = '  -2.40%  ' },
    @{ Row = 24; D = '2.55'; E = '  -1.81%  ' },
    @{ Row = 25; D = $null; E = '  -0.70%  ' },
    @{ Row = 26; D = '1.00'; E = '  +0.17%  ' },
    @{ Row = 27; D = '37.96'; E = '  +0.73%  ' },
    @{ Row = 28; D = '23.16'; E = '  -3.11%  ' },
    @{ Row = 29; D = $null; E = '  +0.53%  ' },
    @{ Row = 30; D = '9.44'; E = '  -0.76%  ' },
    @{ Row = 31; D = '167.07'; E = '  +4.99%  ' },
    @{ Row = 32; D = '0.999'; E = '  -0.08%  ' },
    @{ Row = 33; D = '5.14'; E = '  -2.08%  ' },
    @{ Row = 34; D = $null; E = '  -2.05%  ' },
    @{ Row = 35; D = '17.42'; E = '  +2.26%  ' },
    @{ Row = 36; D = '0.0716'; E = '  -3.28%  ' },
    @{ Row = 37; D = $null; E = '  +0.30%  ' },
    @{ Row = 38; D = $null; E = '  -0.34%  ' },
    @{ Row = 39; D = $null; E = '  -2.91%  ' },
    @{ Row = 40; D = $null; E = '  -2.63%  ' },
    @{ Row = 41; D = $null; E = '  -1.57%  ' },
    @{ Row = 42; D = '1.935.79'; E = '  -3.16%  ' },
    @{ Row = 43; D = '0.0280'; E = '  -1.91%  ' },
    @{ Row = 44; D = '18.52'; E = '  -1.58%  ' },
    @{ Row = 45; D = $null; E = '  -10.57%  ' },
    @{ Row = 46; D = '2.88'; E = '  -2.57%  ' },
    @{ Row = 47; D = '9.62'; E = '  -3.70%  ' },
    @{ Row = 48; D = '53.71'; E = '  +1.53%  ' },
    @{ Row = 49; D = '2.466.02'; E = '  -1.56%  ' },
    @{ Row = 50; D = '70.98'; E = '  -1.62%  ' },
    @{ Row = 51; D = '90.73'; E = '  -1.30%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextCell $ws.Range("D$($u.Row)") $u.D
    }
    $ws.Range("E$($u.Row)").Value = $u.E
}
